$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.019999999999999
$ws.Range("C2").Value = 1.044261553573926
$ws.Range("D2").Value = 1.053165105629876
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.061802322142005
$ws.Range("I2").Value = 1.047876040294763
$ws.Range("J2").Value = 1.049327585015064
$ws.Range("K2").Value = 1.055911828251094
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.064525393959778
$ws.Range("N2").Value = 1.050817750545341
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.045115621805247
$ws.Range("D3").Value = 1.053857968254202
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.06262545716466
$ws.Range("I3").Value = 1.048124720046647
$ws.Range("J3").Value = 1.049829308940285
$ws.Range("K3").Value = 1.056417989191153
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.065163199780368
$ws.Range("N3").Value = 1.051320186976085
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.04566862677661
$ws.Range("D4").Value = 1.054306575992151
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.063158700530994
$ws.Range("I4").Value = 1.04828451081056
$ws.Range("J4").Value = 1.050153651366159
$ws.Range("K4").Value = 1.056745110037366
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.06557586970437
$ws.Range("N4").Value = 1.051644990005406
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.045901196175576
$ws.Range("D5").Value = 1.054495235901522
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.063383022697269
$ws.Range("I5").Value = 1.048351417549271
$ws.Range("J5").Value = 1.050289930392647
$ws.Range("K5").Value = 1.056882534847629
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.065749346752461
$ws.Range("N5").Value = 1.051781462563743
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.045940250608549
$ws.Range("D6").Value = 1.054526916506865
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.063420695924611
$ws.Range("I6").Value = 1.048362635670955
$ws.Range("J6").Value = 1.05031280783692
$ws.Range("K6").Value = 1.056905603365849
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.065778473719021
$ws.Range("N6").Value = 1.051804372496612
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.045671734043757
$ws.Range("D7").Value = 1.054309096620907
$ws.Range("E7").Value = 0.994303590798249
$ws.Range("F7").Value = 1.063161697360723
$ws.Range("I7").Value = 1.048285405880797
$ws.Range("J7").Value = 1.050155472627474
$ws.Range("K7").Value = 1.056746946696296
$ws.Range("L7").Value = 0.9968970624462089
$ws.Range("M7").Value = 1.065578187751905
$ws.Range("N7").Value = 1.051646813853122
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.044550113338886
$ws.Range("D8").Value = 1.053399202979948
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.062080375047668
$ws.Range("I8").Value = 1.047960314566245
$ws.Range("J8").Value = 1.049497207735744
$ws.Range("K8").Value = 1.056082969606591
$ws.Range("L8").Value = 0.9958175282591056
$ws.Range("M8").Value = 1.064740949393644
$ws.Range("N8").Value = 1.050987614149741
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.04257653909855
$ws.Range("D9").Value = 1.051798065809168
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.060179773140955
$ws.Range("I9").Value = 1.047378905534494
$ws.Range("J9").Value = 1.048334964294907
$ws.Range("K9").Value = 1.054909955096418
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.063265439720438
$ws.Range("N9").Value = 1.049823720189906
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.04126283436106
$ws.Range("D10").Value = 1.05073222569792
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.058916054817521
$ws.Range("I10").Value = 1.046985592123648
$ws.Range("J10").Value = 1.047558659418853
$ws.Range("K10").Value = 1.054126000557383
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.062281726624349
$ws.Range("N10").Value = 1.049046312871885
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.040694480559723
$ws.Range("D11").Value = 1.050271101534695
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.058369667307305
$ws.Range("I11").Value = 1.046813938869592
$ws.Range("J11").Value = 1.047222175026753
$ws.Range("K11").Value = 1.053786093685391
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.061855776722761
$ws.Range("N11").Value = 1.048709350633353
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.040483443525009
$ws.Range("D12").Value = 1.05009987997223
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.058166838244783
$ws.Range("I12").Value = 1.046749977733785
$ws.Range("J12").Value = 1.047097139900504
$ws.Range("K12").Value = 1.053659770844537
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.061697562070167
$ws.Range("N12").Value = 1.048584137942884
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.040528708274618
$ws.Range("D13").Value = 1.050136604815235
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.058210340146353
$ws.Range("I13").Value = 1.04676370671979
$ws.Range("J13").Value = 1.047123962596821
$ws.Range("K13").Value = 1.053686870503076
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.061731499531616
$ws.Range("N13").Value = 1.048610998730505
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.040677034637086
$ws.Range("D14").Value = 1.050256947063587
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.058352898862397
$ws.Range("I14").Value = 1.046808655927887
$ws.Range("J14").Value = 1.047211840594426
$ws.Range("K14").Value = 1.053775653145284
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.061842698598864
$ws.Range("N14").Value = 1.048699001524948
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.040768433367819
$ws.Range("D15").Value = 1.050331101958356
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.058440750421549
$ws.Range("I15").Value = 1.046836323943342
$ws.Range("J15").Value = 1.047265978558427
$ws.Range("K15").Value = 1.053830346313215
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.061911212336724
$ws.Range("N15").Value = 1.048753216371067
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.041300564553844
$ws.Range("D16").Value = 1.050762837381716
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.05895233399162
$ws.Range("I16").Value = 1.046996955890243
$ws.Range("J16").Value = 1.047580983728336
$ws.Range("K16").Value = 1.054148549693223
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.062309995732243
$ws.Range("N16").Value = 1.049068668884449
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.041634488302511
$ws.Range("D17").Value = 1.051033759620566
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.059273455185835
$ws.Range("I17").Value = 1.047097356238965
$ws.Range("J17").Value = 1.047778488116902
$ws.Range("K17").Value = 1.054348030815799
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.062560144212663
$ws.Range("N17").Value = 1.049266453751901
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.041829307464596
$ws.Range("D18").Value = 1.051191821656154
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.059460837863821
$ws.Range("I18").Value = 1.047155788149249
$ws.Range("J18").Value = 1.047893656276099
$ws.Range("K18").Value = 1.054464341338856
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.062706051912575
$ws.Range("N18").Value = 1.049381785463093
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.041895743708736
$ws.Range("D19").Value = 1.051245723066271
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.059524743711159
$ws.Range("I19").Value = 1.047175689862404
$ws.Range("J19").Value = 1.047932920042832
$ws.Range("K19").Value = 1.054503992816034
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.062755802678355
$ws.Range("N19").Value = 1.049421104988879
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.041598656563715
$ws.Range("D20").Value = 1.051004688319302
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.059238993827868
$ws.Range("I20").Value = 1.047086597664046
$ws.Range("J20").Value = 1.047757301150618
$ws.Range("K20").Value = 1.054326632863681
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.062533305604153
$ws.Range("N20").Value = 1.049245236697695
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.040633354142132
$ws.Range("D21").Value = 1.050221507588072
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.058310915447952
$ws.Range("I21").Value = 1.046795425062165
$ws.Range("J21").Value = 1.047185964081227
$ws.Range("K21").Value = 1.053749510690249
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.06180995317533
$ws.Range("N21").Value = 1.048673088264132
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.040026864085812
$ws.Range("D22").Value = 1.049729441527967
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.057728111641637
$ws.Range("I22").Value = 1.046611187989156
$ws.Range("J22").Value = 1.046826453983075
$ws.Range("K22").Value = 1.053386268119723
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.061355165697641
$ws.Range("N22").Value = 1.048313067620404
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.040348334342943
$ws.Range("D23").Value = 1.049990261261204
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.058036998518159
$ws.Range("I23").Value = 1.046708965757254
$ws.Range("J23").Value = 1.047017063992137
$ws.Range("K23").Value = 1.05357886577134
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.061596255436903
$ws.Range("N23").Value = 1.048503948317541
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.041614847238765
$ws.Range("D24").Value = 1.051017824272762
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.059254565190996
$ws.Range("I24").Value = 1.047091459402445
$ws.Range("J24").Value = 1.047766874730198
$ws.Range("K24").Value = 1.054336301811539
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.062545432815249
$ws.Range("N24").Value = 1.049254823872855
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.0430864070731
$ws.Range("D25").Value = 1.052211725979618
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.060670541564079
$ws.Range("I25").Value = 1.047530223013258
$ws.Range("J25").Value = 1.048635697214459
$ws.Range("K25").Value = 1.055213555641827
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.063646907744537
$ws.Range("N25").Value = 1.050124880184699
